$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order-line rows to append below the existing data (rows 2-11).
# Every column (including the numeric-looking Quantity/Cost Per/Total Cost
# columns) is stored as TEXT in this workbook, matching the existing rows.
$rows = @(
    @("P4040XC",      "Bag Sheet Pan Cover 30x43",           "1", "31.89", "31.89"),
    @("SAB12032T300",  "Cont Salad - 32oz Sabert (Round)",    "2", "88.05", "176.10"),
    @("4541602",       "Container - Anchor (16oz)",           "1", "43.71", "43.71"),
    @("ANPM424",       "Container - Anchor (24oz)",           "2", "47.17", "94.34"),
    @("ANPLC4LD",      "Lid Anchor - 24/32oz (Dome)",         "2", "56.82", "113.64"),
    @("ANPLC4F",       "Lid Anchor - 24/32oz (Flat)",         "2", "66.89", "133.78"),
    @("SAB52032T300",  "Lid Salad - 24/32oz Sabert (Round)",  "2", "80.81", "161.62"),
    @("PRI80134X60",   "Masking Tape",                        "1", "63.64", "63.64"),
    @("TS12",          "Tamper Evident - 12oz Square",        "1", "38.39", "38.39"),
    @("TS16",          "Tamper Evident - 16oz",                "1", "41.88", "41.88"),
    @("TS8",           "Tamper Evident - 8oz",                 "1", "38.29", "38.29")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Force text storage (no auto-number conversion) on every column - some
    # SKUs (e.g. "4541602") and all of the quantity/cost values look
    # numeric but must stay text, matching the rest of the sheet.
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $data[0]

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $data[1]

    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $data[2]

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $data[3]

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $data[4]
}
